# Apply the "ADDITIONAL SCRAPING" edit:
#  1. Insert a new first sheet "Player Info" with ID/NAME/BATTING_HAND/BOWL_STYLE
#     data for player 6459 (Gudakesh Motie-Kanhai).
#  2. On the "ODI Batting" sheet, rename column D header MATCH_CARD_LINK -> MATCH_CODE
#     and replace the full scorecard URLs with just the numeric match code.
#  3. On the "ODI Bowling" sheet, rename column B header MATCH_CARD_LINK -> MATCH_CODE
#     and replace the full scorecard URLs with just the numeric match code.

$wb = $excel.ActiveWorkbook

# Helper: write a value into $range as TEXT, even when it looks numeric
# (e.g. "6459"), without leaving the cell's number format/style changed -
# mirrors how the source data (scraped, then written as plain strings) looks.
function Set-TextValue {
    param($range, [string]$text)

    $sheet = $range.Worksheet
    $helper = $sheet.Range("ZZ1000")
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy()
    $range.PasteSpecial(-4163) | Out-Null  # xlPasteValues
    $helper.Clear()
}

# ---------------------------------------------------------------------------
# 1. New "Player Info" sheet, inserted before the existing first sheet.
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$playerInfo.Name = "Player Info"

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
$columns = @("A", "B", "C", "D")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $playerInfo.Range("$($columns[$i])1").Value = $headers[$i]
}

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$values = @("6459", "Gudakesh Motie-Kanhai", "Left Handed", "Left Arm Orthodox")
for ($i = 0; $i -lt $values.Length; $i++) {
    Set-TextValue $playerInfo.Range("$($columns[$i])2") $values[$i]
}

# ---------------------------------------------------------------------------
# 2. "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @("4606", "4611", "4616", "4621")
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $row = $i + 2
    Set-TextValue $batting.Range("D$row") $battingCodes[$i]
}

# ---------------------------------------------------------------------------
# 3. "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @("4606", "4611", "4616", "4621")
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $row = $i + 2
    Set-TextValue $bowling.Range("B$row") $bowlingCodes[$i]
}
